# Fruta / hortaliza, semanal
# Insert two new daily-price rows (2021-11-18) for "Vega Monumental Concepción - Limón"
# right before the current row 201, pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 201-202; existing rows 201.. shift down to 203..
$ws.Rows("201:202").Insert()

# --- New row 201 ---
$ws.Cells.Item(201, 1).Value = 11
$ws.Cells.Item(201, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(201, 3).Value = "Bíobío"
$ws.Cells.Item(201, 4).Value = 44518
$ws.Cells.Item(201, 5).Value = 8
$ws.Cells.Item(201, 6).Value = "Fruta"
$ws.Cells.Item(201, 7).Value = 100102
$ws.Cells.Item(201, 8).Value = "Cítricos"
$ws.Cells.Item(201, 9).Value = 100102003
$ws.Cells.Item(201, 10).Value = "Limón"
$ws.Cells.Item(201, 11).Value = "Sin especificar"
$ws.Cells.Item(201, 12).Value = "1a plateado"
$ws.Cells.Item(201, 13).Value = 600
$ws.Cells.Item(201, 14).Value = 8000
$ws.Cells.Item(201, 15).Value = 8500
$ws.Cells.Item(201, 16).Value = 8250
$ws.Cells.Item(201, 17).Value = '$/malla 16 kilos'
$ws.Cells.Item(201, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(201, 19).Value = 516
$ws.Cells.Item(201, 20).Value = 16

# --- New row 202 ---
$ws.Cells.Item(202, 1).Value = 11
$ws.Cells.Item(202, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(202, 3).Value = "Bíobío"
$ws.Cells.Item(202, 4).Value = 44518
$ws.Cells.Item(202, 5).Value = 8
$ws.Cells.Item(202, 6).Value = "Fruta"
$ws.Cells.Item(202, 7).Value = 100102
$ws.Cells.Item(202, 8).Value = "Cítricos"
$ws.Cells.Item(202, 9).Value = 100102003
$ws.Cells.Item(202, 10).Value = "Limón"
$ws.Cells.Item(202, 11).Value = "Sin especificar"
$ws.Cells.Item(202, 12).Value = "2a plateado"
$ws.Cells.Item(202, 13).Value = 300
$ws.Cells.Item(202, 14).Value = 6500
$ws.Cells.Item(202, 15).Value = 6500
$ws.Cells.Item(202, 16).Value = 6500
$ws.Cells.Item(202, 17).Value = '$/malla 16 kilos'
$ws.Cells.Item(202, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(202, 19).Value = 406
$ws.Cells.Item(202, 20).Value = 16
